$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(11, "17-02-2024 05:39", 12, 13, 3, 5),
    @(12, "17-02-2024 05:39", 12, 13, 3, 5),
    @(13, "17-02-2024 05:39", 12, 13, 3, 5),
    @(14, "17-02-2024 05:41", 29, 50, 8, 21),
    @(15, "17-02-2024 05:42", 12, 16, 3, 6),
    @(16, "17-02-2024 05:42", 12, 16, 3, 6),
    @(17, "17-02-2024 05:42", 12, 16, 3, 6),
    @(18, "17-02-2024 05:46", 13, 8, 1, 3),
    @(19, "17-02-2024 05:46", 13, 8, 1, 3),
    @(20, "17-02-2024 05:46", 13, 8, 1, 3),
    @(21, "17-02-2024 05:47", 11, 7, 1, 5),
    @(22, "17-02-2024 05:47", 11, 7, 1, 5),
    @(23, "17-02-2024 05:47", 11, 7, 1, 5),
    @(24, "24-02-2024 14:24", 5, 1, 1, 1),
    @(25, "24-02-2024 14:24", 3, 4, 2, 1),
    @(26, "24-02-2024 14:24", 0, 6, 0, 0),
    @(27, "24-02-2024 14:24", 0, 6, 0, 1),
    @(28, "24-02-2024 14:24", 0, 4, 0, 1),
    @(29, "24-02-2024 14:24", 3, 2, 0, 3),
    @(30, "24-02-2024 14:24", 0, 0, 0, 0),
    @(31, "24-02-2024 19:55", "EMPTY", "EMPTY", "EMPTY", "EMPTY"),
    @(32, "24-02-2024 19:58", 11, 24, 4, 7),
    @(33, "26-02-2024 20:41", 12, 13, 2, 7)
)

foreach ($row in $data) {
    $r = $row[0]
    $ts = $row[1]
    $ws.Cells.Item($r, 1).Value = $ts
    for ($col = 2; $col -le 5; $col++) {
        $v = $row[$col]
        if ($v -ne "EMPTY") {
            $ws.Cells.Item($r, $col).Value = $v
        }
    }
}
